$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.104.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.549.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3808'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3268'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.126'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07325'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.788'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.752'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.550.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001082'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06586'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.380'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.094.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.290'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.491'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.910'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.83%  '
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.722.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.850'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.863'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08204'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.222'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.85%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02308'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.34%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06223'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.237'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2147'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.236'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5990'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.00%  '
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.723'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5788'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.967'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.170'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.35%  '
$ws.Range("E51").Value = '  -2.94%  '
